# Auto-generated script to apply scheduled market-data refresh values
# to the Chocobo_Profits workbook (H/I/J/K/L/M/N price & profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 279.625
$ws.Range("I11").Value = 279.625
$ws.Range("K11").Value = 279.625
$ws.Range("M11").Value = -139.625

$ws.Range("H19").Value = 1254188.4
$ws.Range("J19").Value = 1707.5
$ws.Range("L19").Value = 1707.5
$ws.Range("N19").Value = -2057.5

$ws.Range("H40").Value = 1714.4445
$ws.Range("I40").Value = 1587.9
$ws.Range("J40").Value = 1872.625
$ws.Range("K40").Value = 1587.9
$ws.Range("L40").Value = 1872.625
$ws.Range("M40").Value = -1412.9
$ws.Range("N40").Value = -2222.625

$ws.Range("H41").Value = 524.6923
$ws.Range("I41").Value = 145.57143
$ws.Range("J41").Value = 967
$ws.Range("K41").Value = 145.57143
$ws.Range("L41").Value = 967
$ws.Range("M41").Value = 294.42857
$ws.Range("N41").Value = -1847

$ws.Range("H86").Value = 1427.375
$ws.Range("I86").Value = 1352.8334
$ws.Range("J86").Value = 1651
$ws.Range("K86").Value = 1352.8334
$ws.Range("L86").Value = 1651
$ws.Range("M86").Value = -229.8334
$ws.Range("N86").Value = -3897

$ws.Range("H89").Value = 1427.375
$ws.Range("I89").Value = 1352.8334
$ws.Range("J89").Value = 1651
$ws.Range("K89").Value = 6764.166999999999
$ws.Range("L89").Value = 8255
$ws.Range("M89").Value = -1148.166999999999
$ws.Range("N89").Value = -19487

$ws.Range("H107").Value = 1188.3478
$ws.Range("I107").Value = 1160.8125
$ws.Range("J107").Value = 1251.2858
$ws.Range("K107").Value = 1160.8125
$ws.Range("L107").Value = 1251.2858
$ws.Range("M107").Value = 759.1875
$ws.Range("N107").Value = -5091.2858

$ws.Range("H113").Value = 10659.071
$ws.Range("I113").Value = 3188.75
$ws.Range("J113").Value = 13647.2
$ws.Range("K113").Value = 3188.75
$ws.Range("L113").Value = 13647.2
$ws.Range("M113").Value = 65.25
$ws.Range("N113").Value = -20155.2


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 313
$ws.Range("I5").Value = 141.25
$ws.Range("K5").Value = 141.25
$ws.Range("M5").Value = -29.25

$ws.Range("H88").Value = 4169622
$ws.Range("I88").Value = 13334933
$ws.Range("J88").Value = 3571.2727
$ws.Range("K88").Value = 13334933
$ws.Range("L88").Value = 3571.2727
$ws.Range("M88").Value = -13334527
$ws.Range("N88").Value = -4383.2727

$ws.Range("H91").Value = 4169622
$ws.Range("I91").Value = 13334933
$ws.Range("J91").Value = 3571.2727
$ws.Range("K91").Value = 13334933
$ws.Range("L91").Value = 3571.2727
$ws.Range("M91").Value = -13333529
$ws.Range("N91").Value = -6379.2727

$ws.Range("H92").Value = 26275
$ws.Range("J92").Value = 26275
$ws.Range("L92").Value = 26275
$ws.Range("N92").Value = -31267


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 313
$ws.Range("I4").Value = 141.25
$ws.Range("K4").Value = 141.25
$ws.Range("M4").Value = -26.25

$ws.Range("H38").Value = 16994
$ws.Range("J38").Value = 16994
$ws.Range("L38").Value = 16994
$ws.Range("N38").Value = -17826

$ws.Range("H86").Value = 2028.7368
$ws.Range("I86").Value = 1648.1538
$ws.Range("J86").Value = 2853.3333
$ws.Range("K86").Value = 1648.1538
$ws.Range("L86").Value = 2853.3333
$ws.Range("M86").Value = -525.1538
$ws.Range("N86").Value = -5099.3333

$ws.Range("H88").Value = 48600
$ws.Range("J88").Value = 48600
$ws.Range("L88").Value = 48600
$ws.Range("N88").Value = -49412

$ws.Range("H89").Value = 2028.7368
$ws.Range("I89").Value = 1648.1538
$ws.Range("J89").Value = 2853.3333
$ws.Range("K89").Value = 8240.769
$ws.Range("L89").Value = 14266.6665
$ws.Range("M89").Value = -2624.769
$ws.Range("N89").Value = -25498.6665

$ws.Range("H91").Value = 48600
$ws.Range("J91").Value = 48600
$ws.Range("L91").Value = 48600
$ws.Range("N91").Value = -51408

$ws.Range("H107").Value = 1801.909
$ws.Range("I107").Value = 1644.4
$ws.Range("J107").Value = 1933.1666
$ws.Range("K107").Value = 1644.4
$ws.Range("L107").Value = 1933.1666
$ws.Range("M107").Value = 275.5999999999999
$ws.Range("N107").Value = -5773.1666

$ws.Range("H123").Value = 29885
$ws.Range("J123").Value = 29885
$ws.Range("L123").Value = 29885
$ws.Range("N123").Value = -39685


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12347709
$ws.Range("I16").Value = 27779248
$ws.Range("J16").Value = 2477.2
$ws.Range("K16").Value = 27779248
$ws.Range("L16").Value = 2477.2
$ws.Range("M16").Value = -27778961
$ws.Range("N16").Value = -3051.2

$ws.Range("H31").Value = 4766.377
$ws.Range("I31").Value = 1821.3334
$ws.Range("J31").Value = 10288.333
$ws.Range("K31").Value = 1821.3334
$ws.Range("L31").Value = 10288.333
$ws.Range("M31").Value = -1526.3334
$ws.Range("N31").Value = -10878.333

$ws.Range("H34").Value = 4766.377
$ws.Range("I34").Value = 1821.3334
$ws.Range("J34").Value = 10288.333
$ws.Range("K34").Value = 1821.3334
$ws.Range("L34").Value = 10288.333
$ws.Range("M34").Value = -1619.3334
$ws.Range("N34").Value = -10692.333

$ws.Range("H58").Value = 2204.551
$ws.Range("I58").Value = 1378.4722
$ws.Range("J58").Value = 4492.154
$ws.Range("K58").Value = 1378.4722
$ws.Range("L58").Value = 4492.154
$ws.Range("M58").Value = -1175.4722
$ws.Range("N58").Value = -4898.154

$ws.Range("H62").Value = 2742.2
$ws.Range("I62").Value = 2826.25
$ws.Range("J62").Value = 2406
$ws.Range("K62").Value = 2826.25
$ws.Range("L62").Value = 2406
$ws.Range("M62").Value = -2202.25
$ws.Range("N62").Value = -3654

$ws.Range("H65").Value = 2742.2
$ws.Range("I65").Value = 2826.25
$ws.Range("J65").Value = 2406
$ws.Range("K65").Value = 14131.25
$ws.Range("L65").Value = 12030
$ws.Range("M65").Value = -11011.25
$ws.Range("N65").Value = -18270

$ws.Range("H88").Value = 31827.75
$ws.Range("J88").Value = 39000
$ws.Range("L88").Value = 39000
$ws.Range("N88").Value = -39812

$ws.Range("H91").Value = 31827.75
$ws.Range("J91").Value = 39000
$ws.Range("L91").Value = 39000
$ws.Range("N91").Value = -41808

$ws.Range("H99").Value = 8336736
$ws.Range("I99").Value = 12501418
$ws.Range("K99").Value = 12501418
$ws.Range("M99").Value = -12499920

$ws.Range("H113").Value = 12347709
$ws.Range("I113").Value = 27779248
$ws.Range("J113").Value = 2477.2
$ws.Range("K113").Value = 27779248
$ws.Range("L113").Value = 2477.2
$ws.Range("M113").Value = -27777078
$ws.Range("N113").Value = -6817.2

$ws.Range("H122").Value = 2223.1428
$ws.Range("I122").Value = 1390.5
$ws.Range("J122").Value = 3333.3333
$ws.Range("K122").Value = 4171.5
$ws.Range("L122").Value = 9999.999899999999
$ws.Range("M122").Value = -1721.5
$ws.Range("N122").Value = -14899.9999

$ws.Range("H126").Value = 8336736
$ws.Range("I126").Value = 12501418
$ws.Range("K126").Value = 37504254
$ws.Range("M126").Value = -37501784

$ws.Range("H136").Value = 2204.551
$ws.Range("I136").Value = 1378.4722
$ws.Range("J136").Value = 4492.154
$ws.Range("K136").Value = 4135.4166
$ws.Range("L136").Value = 13476.462
$ws.Range("M136").Value = -1585.4166
$ws.Range("N136").Value = -18576.462


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M62").ClearContents()
$ws.Range("H62").Value = 9000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 27000
$ws.Range("N62").Value = -28372

$ws.Range("M65").ClearContents()
$ws.Range("H65").Value = 9000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 81000
$ws.Range("N65").Value = -87864

$ws.Range("H113").Value = 5208903
$ws.Range("I113").Value = 617.9091
$ws.Range("K113").Value = 1853.7273
$ws.Range("M113").Value = 316.2727

$ws.Range("H131").Value = 777.99
$ws.Range("J131").Value = 829.9888999999999
$ws.Range("L131").Value = 2489.9667
$ws.Range("N131").Value = -12569.9667


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 182.47058
$ws.Range("I2").Value = 61
$ws.Range("J2").Value = 356
$ws.Range("K2").Value = 61
$ws.Range("L2").Value = 356
$ws.Range("M2").Value = 52
$ws.Range("N2").Value = -582

$ws.Range("H70").Value = 6173.8374
$ws.Range("I70").Value = 5831.2856
$ws.Range("J70").Value = 6813.2666
$ws.Range("K70").Value = 5831.2856
$ws.Range("L70").Value = 6813.2666
$ws.Range("M70").Value = -5561.2856
$ws.Range("N70").Value = -7353.2666

$ws.Range("H73").Value = 6173.8374
$ws.Range("I73").Value = 5831.2856
$ws.Range("J73").Value = 6813.2666
$ws.Range("K73").Value = 5831.2856
$ws.Range("L73").Value = 6813.2666
$ws.Range("M73").Value = -4895.2856
$ws.Range("N73").Value = -8685.266599999999

$ws.Range("H102").Value = 1909.025
$ws.Range("I102").Value = 1404.2693
$ws.Range("J102").Value = 2846.4285
$ws.Range("K102").Value = 1404.2693
$ws.Range("L102").Value = 2846.4285
$ws.Range("M102").Value = 217.7307000000001
$ws.Range("N102").Value = -6090.4285

$ws.Range("N138").ClearContents()
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0

$ws.Range("H140").Value = 40731.11
$ws.Range("J140").Value = 40731.11
$ws.Range("L140").Value = 40731.11
$ws.Range("N140").Value = -51091.11


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3537.6316
$ws.Range("I132").Value = 2815.1035
$ws.Range("J132").Value = 5865.778
$ws.Range("K132").Value = 8445.3105
$ws.Range("L132").Value = 17597.334
$ws.Range("M132").Value = -5915.3105
$ws.Range("N132").Value = -22657.334

$ws.Range("H141").Value = 32247.5
$ws.Range("J141").Value = 32247.5
$ws.Range("L141").Value = 32247.5
$ws.Range("N141").Value = -42607.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 39466.332
$ws.Range("J80").Value = 39466.332
$ws.Range("L80").Value = 39466.332
$ws.Range("N80").Value = -41462.332

$ws.Range("H83").Value = 39466.332
$ws.Range("J83").Value = 39466.332
$ws.Range("L83").Value = 118398.996
$ws.Range("N83").Value = -128382.996

$ws.Range("H126").Value = 2085.7666
$ws.Range("I126").Value = 963.2727
$ws.Range("K126").Value = 2889.8181
$ws.Range("M126").Value = -419.8181

